$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[V.%Tsipropoulou%NULL%1,                           V. A.%Nikopoulou%NULL%1,                           Vasiliki%Holeva%vholeva@yahoo.gr%0,                           Z.%Nasika%NULL%2,                           Z.%Nasika%NULL%0,                           I.%Diakogiannis%NULL%1,                           S.%Sakka%NULL%1,                           S.%Kostikidou%NULL%1,                           C.%Varvara%NULL%1,                           E.%Spyridopoulou%NULL%1,                           E.%Parlapani%NULL%1]"
$ws.Range("E3").Value = "[Catherine A.%Perz%perzc@uhv.edu%1,                           Brent A.%Lang%NULL%2,                           Brent A.%Lang%NULL%0,                           Rick%Harrington%NULL%1]"
$ws.Range("E4").Value = "[Kun-Chia%Chang%NULL%1,                           Wen-Li%Hou%NULL%1,                           Amir H.%Pakpour%NULL%4,                           Chung-Ying%Lin%cylin36933@gmail.com%3,                           Chung-Ying%Lin%cylin36933@gmail.com%0,                           Mark D.%Griffiths%NULL%5,                           Mark D.%Griffiths%NULL%0]"
$ws.Range("E5").Value = "[Lin-sen%Feng%NULL%1,                           Zheng-jiao%Dong%NULL%1,                           Ruo-yu%Yan%NULL%1,                           Xiao-qian%Wu%NULL%1,                           Li%Zhang%NULL%0,                           Jun%Ma%NULL%1,                           Yong%Zeng%NULL%1]"
$ws.Range("E6").Value = "[Abdulkadir%Haktanir%xref no email%1,       Tolga%Seki%xref no email%1,       B\u00fclent%Dilma\u00e7%xref no email%1]"
$ws.Range("E7").Value = "[Daniel Kwasi%Ahorsu%NULL%1,                           Chung-Ying%Lin%NULL%0,                           Vida%Imani%NULL%1,                           Mohsen%Saffari%NULL%1,                           Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%0,                           Amir H.%Pakpour%pakpour_amir@yahoo.com%0]"
$ws.Range("E8").Value = "[Mohsen%Alyami%maly002@aucklanduni.ac.nz%1,                           Marcus%Henning%NULL%2,                           Marcus%Henning%NULL%0,                           Christian U.%Krägeloh%NULL%2,                           Christian U.%Krägeloh%NULL%0,                           Hussain%Alyami%NULL%2,                           Hussain%Alyami%NULL%0]"
$ws.Range("E9").Value = "[Dana%Tzur Bitan%NULL%0,                           Ariella%Grossman-Giron%NULL%1,                           Yuval%Bloch%NULL%1,                           Yael%Mayer%NULL%1,                           Noga%Shiffman%NULL%1,                           Shlomo%Mendlovic%NULL%1]"
$ws.Range("E10").Value = "[Najmuj%Sakib%NULL%1,                           A. K. M. Israfil%Bhuiyan%NULL%1,                           Sahadat%Hossain%NULL%1,                           Firoj%Al Mamun%NULL%1,                           Ismail%Hosen%NULL%1,                           Abu Hasnat%Abdullah%NULL%1,                           Md. Abedin%Sarker%NULL%1,                           Mohammad Sarif%Mohiuddin%NULL%1,                           Istihak%Rayhan%NULL%1,                           Moazzem%Hossain%NULL%1,                           Md. Tajuddin%Sikder%NULL%1,                           David%Gozal%NULL%1,                           Mohammad%Muhit%NULL%1,                           S. M. Shariful%Islam%NULL%1,                           Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%0,                           Amir H.%Pakpour%apakpour@qums.ac.ir%0,                           Amir H.%Pakpour%apakpour@qums.ac.ir%0,                           Mohammed A.%Mamun%mamunphi46@gmail.com%2,                           Mohammed A.%Mamun%mamunphi46@gmail.com%0]"
$ws.Range("E11").Value = "[Paolo%Soraci%paolo.soraci85@gmail.com%0,                           Ambra%Ferrari%a.ferrari50@campus.unimib.it%1,                           Francesco A.%Abbiati%francesco.abbiati@edu.unito.it%1,                           Elena%Del Fante%edelfante@hotmail.it%1,                           Rosanna%De Pace%Rosannadepace@libero.it%1,                           Antonino%Urso%antonino.urso@tin.it%1,                           Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%0]"
